$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update pizza ingredient inventory counts in row 2
$ws.Range("B2:D2").Value = 990
$ws.Range("E2:F2").Value = 995
